$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ct manual (1-5)")
$ws.Rows.Item(3).Delete()
